$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.429807782173157
$ws.Range("B1").Value = 2.734714508056641
$ws.Range("C1").Value = 6.364161491394043
$ws.Range("D1").Value = 2.351888179779053
$ws.Range("E1").Value = 1.162959575653076
